# Add CAAF Data backtest
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix existing row 7 (Lambda/Epsilon column D) value
$ws.Range("D7").Value = 0.025

# New rows 8-10
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "JP"
$ws.Range("C8").Value = 0.07
$ws.Range("D8").Value = 0.025
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = "None"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "US"
$ws.Range("C9").Value = 0.07
$ws.Range("D9").Value = 0.025
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = "None"

$ws.Range("A10").Value = 9
$ws.Range("C10").Value = 0.07
$ws.Range("D10").Value = 0.025
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = "Yes"

# Selection / window view tweaks
$ws.Range("Q15").Select() | Out-Null
